# Refresh the cryptos list with latest prices / 1h volume changes.
# Note: some Price values (column D) look like plain decimal numbers
# (e.g. "0.4830", "1.100"); a leading apostrophe is used so Excel stores
# them as literal text (preserving trailing zeros) instead of coercing
# them to numeric values, matching how the sheet already stores these
# cells as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.956.07'
$ws.Range('E2').Value = '  +1.83%  '
$ws.Range('D3').Value = '1.904.46'
$ws.Range('E3').Value = '  +2.48%  '
$ws.Range('D4').Value = '''1.006'
$ws.Range('D5').Value = '''317.13'
$ws.Range('E5').Value = '  +1.70%  '
$ws.Range('E6').Value = '  -0.48%  '
$ws.Range('D7').Value = '''0.4830'
$ws.Range('E7').Value = '  +1.27%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '''0.07373'
$ws.Range('E9').Value = '  +0.93%  '
$ws.Range('D10').Value = '''0.9327'
$ws.Range('E10').Value = '  +0.41%  '
$ws.Range('D11').Value = '''20.77'
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('D12').Value = '''0.07735'
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('D13').Value = '1.943.45'
$ws.Range('E13').Value = '  +4.07%  '
$ws.Range('D14').Value = '''5.485'
$ws.Range('E14').Value = '  +0.67%  '
$ws.Range('D15').Value = '''6.624'
$ws.Range('E15').Value = '  +1.11%  '
$ws.Range('D16').Value = '''91.69'
$ws.Range('E16').Value = '  +1.71%  '
$ws.Range('D17').Value = '''1.006'
$ws.Range('E17').Value = '  -0.58%  '
$ws.Range('D18').Value = '''0.000008866'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('E19').Value = '  -0.43%  '
$ws.Range('D20').Value = '28.009.35'
$ws.Range('E20').Value = '  +1.93%  '
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('D22').Value = '''5.147'
$ws.Range('E22').Value = '  +1.11%  '
$ws.Range('D23').Value = '2.175.29'
$ws.Range('E23').Value = '  +3.60%  '
$ws.Range('D24').Value = '''10.89'
$ws.Range('E24').Value = '  +1.96%  '
$ws.Range('D25').Value = '''156.12'
$ws.Range('E25').Value = '  +0.83%  '
$ws.Range('D26').Value = '''1.913'
$ws.Range('E26').Value = '  -1.38%  '
$ws.Range('E27').Value = '  +0.23%  '
$ws.Range('D28').Value = '''2.122'
$ws.Range('E28').Value = '  +6.13%  '
$ws.Range('D29').Value = '''117.31'
$ws.Range('E29').Value = '  +1.87%  '
$ws.Range('D30').Value = '''4.967'
$ws.Range('E30').Value = '  +0.77%  '
$ws.Range('D31').Value = '''0.08948'
$ws.Range('E31').Value = '  +0.66%  '
$ws.Range('D32').Value = '''3.246'
$ws.Range('E32').Value = '  -2.56%  '
$ws.Range('D33').Value = '''1.257'
$ws.Range('E33').Value = '  +4.58%  '
$ws.Range('D34').Value = '''0.7663'
$ws.Range('E34').Value = '  +2.16%  '
$ws.Range('D35').Value = '''4.668'
$ws.Range('E35').Value = '  +2.09%  '
$ws.Range('E36').Value = '  +0.70%  '
$ws.Range('D37').Value = '''2.548'
$ws.Range('E37').Value = '  -5.91%  '
$ws.Range('D38').Value = '''1.100'
$ws.Range('E38').Value = '  -1.91%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '''0.05276'
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '''0.5482'
$ws.Range('E40').Value = '  -1.38%  '
$ws.Range('D41').Value = '''2.998'
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('D42').Value = '''6.948'
$ws.Range('E42').Value = '  -0.99%  '
$ws.Range('E43').Value = '  +0.65%  '
$ws.Range('D44').Value = '''8.478'
$ws.Range('E44').Value = '  -1.03%  '
$ws.Range('D45').Value = '''110.09'
$ws.Range('E45').Value = '  +6.74%  '
$ws.Range('D46').Value = '''10.67'
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').Value = '''0.4803'
$ws.Range('E47').Value = '  -1.21%  '
$ws.Range('E48').Value = '  -0.49%  '
$ws.Range('D49').Value = '''1.647'
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('D50').Value = '''67.82'
$ws.Range('E50').Value = '  +0.82%  '
$ws.Range('D51').Value = '''0.06080'
